$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '22.443.32'
Set-TextValue $ws.Range('E2') '  -4.36%  '
Set-TextValue $ws.Range('D3') '1.572.31'
Set-TextValue $ws.Range('E3') '  -4.64%  '
Set-TextValue $ws.Range('D4') '0.9995'
Set-TextValue $ws.Range('E4') '  -0.08%  '
Set-TextValue $ws.Range('D5') '1.000'
Set-TextValue $ws.Range('E5') '  -0.02%  '
Set-TextValue $ws.Range('D6') '291.24'
Set-TextValue $ws.Range('D7') '0.3656'
Set-TextValue $ws.Range('E7') '  -3.36%  '
Set-TextValue $ws.Range('D8') '49.27'
Set-TextValue $ws.Range('E8') '  -1.24%  '
Set-TextValue $ws.Range('D9') '0.3380'
Set-TextValue $ws.Range('E9') '  -4.96%  '
Set-TextValue $ws.Range('E10') '  -3.91%  '
Set-TextValue $ws.Range('D11') '0.07592'
Set-TextValue $ws.Range('E11') '  -6.22%  '
Set-TextValue $ws.Range('D12') '1.000'
Set-TextValue $ws.Range('E12') '  +0.00%  '
Set-TextValue $ws.Range('E13') '  -4.17%  '
Set-TextValue $ws.Range('D14') '6.064'
Set-TextValue $ws.Range('E14') '  -5.21%  '
Set-TextValue $ws.Range('D15') '6.882'
Set-TextValue $ws.Range('E15') '  -6.48%  '
Set-TextValue $ws.Range('D16') '0.00001142'
Set-TextValue $ws.Range('E16') '  -4.60%  '
Set-TextValue $ws.Range('D17') '1.569.09'
Set-TextValue $ws.Range('E17') '  -5.08%  '
Set-TextValue $ws.Range('D18') '89.13'
Set-TextValue $ws.Range('E18') '  -8.40%  '
Set-TextValue $ws.Range('D19') '0.06730'
Set-TextValue $ws.Range('E19') '  -3.07%  '
Set-TextValue $ws.Range('E20') '  +0.06%  '
Set-TextValue $ws.Range('D21') '6.274'
Set-TextValue $ws.Range('E21') '  -7.10%  '
Set-TextValue $ws.Range('D22') '16.47'
Set-TextValue $ws.Range('E22') '  -4.94%  '
Set-TextValue $ws.Range('E23') '  -9.00%  '
Set-TextValue $ws.Range('D24') '12.04'
Set-TextValue $ws.Range('E24') '  -3.32%  '
Set-TextValue $ws.Range('D25') '22.448.80'
Set-TextValue $ws.Range('E25') '  -4.34%  '
Set-TextValue $ws.Range('D26') '2.391'
Set-TextValue $ws.Range('E26') '  -4.11%  '
Set-TextValue $ws.Range('D27') '3.012'
Set-TextValue $ws.Range('E27') '  +3.64%  '
Set-TextValue $ws.Range('D28') '19.90'
Set-TextValue $ws.Range('E28') '  -4.87%  '
Set-TextValue $ws.Range('D29') '144.19'
Set-TextValue $ws.Range('E29') '  -5.82%  '
Set-TextValue $ws.Range('D30') '5.003'
Set-TextValue $ws.Range('E30') '  -3.81%  '
Set-TextValue $ws.Range('D31') '125.18'
Set-TextValue $ws.Range('E31') '  -5.77%  '
Set-TextValue $ws.Range('D32') '1.747.13'
Set-TextValue $ws.Range('E32') '  -4.68%  '
Set-TextValue $ws.Range('D33') '1.041'
Set-TextValue $ws.Range('E33') '  +4.55%  '
Set-TextValue $ws.Range('D34') '6.301'
Set-TextValue $ws.Range('E34') '  -9.10%  '
Set-TextValue $ws.Range('D35') '1.971'
Set-TextValue $ws.Range('E35') '  -7.29%  '
Set-TextValue $ws.Range('D36') '10.41'
Set-TextValue $ws.Range('E36') '  -8.96%  '
Set-TextValue $ws.Range('D37') '0.02563'
Set-TextValue $ws.Range('E37') '  -5.73%  '
Set-TextValue $ws.Range('D38') '0.08442'
Set-TextValue $ws.Range('E38') '  -3.37%  '
Set-TextValue $ws.Range('D39') '0.2306'
Set-TextValue $ws.Range('E39') '  -5.19%  '
Set-TextValue $ws.Range('D40') '0.06547'
Set-TextValue $ws.Range('E40') '  -3.45%  '
Set-TextValue $ws.Range('D41') '5.517'
Set-TextValue $ws.Range('E41') '  -7.10%  '
Set-TextValue $ws.Range('D42') '11.88'
Set-TextValue $ws.Range('E42') '  -9.17%  '
Set-TextValue $ws.Range('D43') '1.255'
Set-TextValue $ws.Range('E43') '  -3.58%  '
Set-TextValue $ws.Range('D44') '0.6403'
Set-TextValue $ws.Range('E44') '  -7.11%  '
Set-TextValue $ws.Range('E45') '  -6.81%  '
Set-TextValue $ws.Range('D46') '0.9994'
Set-TextValue $ws.Range('E46') '  -0.03%  '
Set-TextValue $ws.Range('D47') '0.6038'
Set-TextValue $ws.Range('E47') '  -5.10%  '
Set-TextValue $ws.Range('D48') '3.776'
Set-TextValue $ws.Range('E48') '  -3.36%  '
Set-TextValue $ws.Range('D49') '2.140'
Set-TextValue $ws.Range('E49') '  -5.25%  '
Set-TextValue $ws.Range('D50') '122.67'
Set-TextValue $ws.Range('E50') '  -3.56%  '
Set-TextValue $ws.Range('D51') '1.207'
Set-TextValue $ws.Range('E51') '  +2.17%  '
